$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.016143202781677
$ws.Range("B1").Value = 1.142959356307983
$ws.Range("C1").Value = 5.500528812408447
$ws.Range("D1").Value = 1.626428484916687
$ws.Range("E1").Value = 0.9922494888305664
